$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = "#08-파이썬(Python) 컴프리헨션(Comprehension)"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-08"

# Row 9
$ws.Range("D9").Value = "[공지] MSDS 입학시험 문제 공개"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msds-prep-exam-202106/#utm_source=rss&utm_medium=rss&utm_campaign=msds-prep-exam-202106"

# Row 41
$ws.Range("D41").Value = "ML 모델 도입을 위한 SageMaker의 효율성"
$ws.Range("E41").Value = "http://cloudinsight.net/ai/%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d-%eb%aa%a8%eb%8d%b8-%ea%b0%9c%eb%b0%9c%ea%b3%bc-sagemaker/"

# Row 51
$ws.Range("D51").Value = "[github] 깃헙에 게시한 이슈 삭제하려면"
$ws.Range("E51").Value = "https://bskyvision.com/1205"
